# "Transcription and Stylesheet Work"
#
# Appends two new glyph entries to the "Glyphs" sheet's code list:
#   g54 -> colon
#   g55 -> da with stroke
# (continuing the existing g0..g53 / description table in columns A:B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A55").Value = "g54"
$ws.Range("B55").Value = "colon"
$ws.Range("A56").Value = "g55"
$ws.Range("B56").Value = "da with stroke"

# Match the saved workbook's final on-screen focus: the newly added last row.
[void]$ws.Range("B56").Select()
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
